$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.120.86"
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").Value = "1.562.70"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'206.21"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "'0.493"
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'22.11"
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "1.785.25"
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").Value = "1.549.11"
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("E15").Value = "  -2.79%  "
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "27.130.19"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").Value = "0.0₃0687"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").Value = "'211.82"
$ws.Range("E20").Value = "  -1.71%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'152.14"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").Value = "'6.59"
$ws.Range("E26").Value = "  -4.11%  "
$ws.Range("D27").Value = "'14.86"
$ws.Range("E27").Value = "  -1.95%  "
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "'3.17"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").Value = "1.377.06"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("D39").Value = "'0.521"
$ws.Range("E39").Value = "  -3.21%  "
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'0.994"
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("D43").Value = "'1.77"
$ws.Range("E43").Value = "  +2.51%  "
$ws.Range("D44").Value = "'63.44"
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").Value = "'5.20"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").Value = "1.696.77"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").Value = "'85.38"
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("E51").Value = "  +0.07%  "
# Row 36/37: HuobiToken and TrustWalletToken swap positions
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'0.943"
$ws.Range("E36").Value = "  -3.92%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.28"
$ws.Range("E37").Value = "  -1.22%  "